$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21, pushing the existing weekly records (old rows 21-33)
# down to rows 22-34.
$ws.Rows.Item(21).Insert()

# Populate the new weekly record in row 21.
$ws.Range("A21").Value = 3
$ws.Range("B21").Value = "Femacal de La Calera"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44452
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 100112035
$ws.Range("G21").Value = "Bruselas (repollito)"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 73
$ws.Range("K21").Value = 22000
$ws.Range("L21").Value = 23000
$ws.Range("M21").Value = 22479
$ws.Range("N21").Value = "$/malla 15 kilos"
$ws.Range("O21").Value = "Provincia de Quillota"
$ws.Range("P21").Value = 1499
$ws.Range("Q21").Value = 15
$ws.Range("R21").Value = "Hortaliza"
